$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config_Spec")
$tbl = $ws.ListObjects.Item(1)

# Update Flag in row 11 from Y to N
$ws.Range("B11").Value = "N"

# Add new row to the table (auto-expands table range)
$newRow = $tbl.ListRows.Add()

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Y"
$ws.Range("C12").Value = "LoginScreen.xlsx"
$ws.Range("D12").Value = "LoginScreen"
$ws.Range("E12").Value = "LoginScreen.xlsx"
$ws.Range("F12").Value = "LoginScreen"
$ws.Range("G12").Value = "Y"
$ws.Range("H12").Value = "Login"
$ws.Range("I12").Value = "LoginScreen"
$ws.Range("J12").Value = "LoginScreen"

# Copy formatting from row 11 down to row 12 (after values so numeric stays numeric)
$ws.Range("A11:J11").Copy()
$ws.Range("A12:J12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Activate()
$ws.Range("A13").Select() | Out-Null

Write-Output "done"
